# Applies the textual corrections / updates described in the commit:
# "From 1.2.4 to 1.2.5 change and minor updates"
#
# The underlying change set (derived from the shared-strings diff) amounts to:
#  - bump the displayed Version from 1.0 to 1.2.5
#  - fix a typo ("usuario" -> "usuário") and add a trailing period to the
#    Precondition text (repeated for every test case block)
#  - add trailing periods to two recurring step texts
#  - fix "conta bancários" -> "conta bancária" in the TC2 MSG403 text
#  - swap the TC3 "Confirma" expected result so it now shows the MSG213
#    communication-failure message, while TC1's "Confirma" expected result
#    now shows the success message that used to live on TC3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version bump (header area)
$ws.Range("D2").Value = "1.2.5"

# Precondition text: fix typo + add trailing period (same text appears on
# each of the four test-case blocks)
$precondition = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B8").Value  = $precondition
$ws.Range("B18").Value = $precondition
$ws.Range("B27").Value = $precondition
$ws.Range("B37").Value = $precondition

# "Chefe Acessa..." step: add trailing period (same text on each block)
$acessa = "Chefe Acessa a funcionalidade Minha Conta Bancária (menu)."
$ws.Range("B10").Value = $acessa
$ws.Range("B20").Value = $acessa
$ws.Range("B29").Value = $acessa
$ws.Range("B39").Value = $acessa

# "SYSTEM Apresenta os campos..." step: add trailing period
$apresenta = "SYSTEM Apresenta os campos (banco/agência/conta corrente) alterados."
$ws.Range("D11").Value = $apresenta
$ws.Range("D30").Value = $apresenta
$ws.Range("D40").Value = $apresenta

# TC2 MSG403 text: "conta bancários" -> "conta bancária"
$ws.Range("D21").Value = "SYSTEM Exibe mensagens informativas (MSG403 - Informativos sobre a atualização de conta bancária (dados bancários)) para o usuário sobre a manutenção de informações bancárias."

# TC1 "Confirma" expected result now holds the success message ...
$ws.Range("D13").Value = "SYSTEM Atualiza os dados bancários do beneficiário na base do RH (SRH); Exibe mensagem de sucesso para o usuário."

# ... and TC3 "Confirma" expected result now holds the MSG213 failure message
$ws.Range("D32").Value = "SYSTEM Identifica que ocorreu uma falha durante a tentativa de atualização dos dados bancários; Mantém os dados consistentes, interrompe a operação; Exibe mensagem de erro (MSG213 - Não foi possível concluir a operação. Falha na comunicação com o sistema de Recursos Humanos) para o usuário."
